$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 135.5625
$ws.Range("I33").Value = 117.23077
$ws.Range("J33").Value = 215
$ws.Range("K33").Value = 117.23077
$ws.Range("L33").Value = 215
$ws.Range("M33").Value = 111.76923
$ws.Range("N33").Value = -673

$ws.Range("H45").Value = 391.66666
$ws.Range("I45").Value = 275
$ws.Range("K45").Value = 825
$ws.Range("M45").Value = -633

$ws.Range("H64").Value = 9805.454
$ws.Range("I64").Value = 3965
$ws.Range("K64").Value = 3965
$ws.Range("M64").Value = -3717

$ws.Range("H67").Value = 9805.454
$ws.Range("I67").Value = 3965
$ws.Range("K67").Value = 3965
$ws.Range("M67").Value = -3107

$ws.Range("H113").Value = 250002220
$ws.Range("I113").Value = 1000000000
$ws.Range("J113").Value = 2966.6667
$ws.Range("K113").Value = 1000000000
$ws.Range("L113").Value = 2966.6667
$ws.Range("M113").Value = -999996746
$ws.Range("N113").Value = -9474.6667

$ws.Range("H116").Value = 3723.2173
$ws.Range("I116").Value = 3270.5
$ws.Range("J116").Value = 4427.4443
$ws.Range("K116").Value = 3270.5
$ws.Range("L116").Value = 4427.4443
$ws.Range("M116").Value = 171.5
$ws.Range("N116").Value = -11311.4443

$ws.Range("H132").Value = 8554.143
$ws.Range("I132").Value = 1563.1666
$ws.Range("J132").Value = 50500
$ws.Range("K132").Value = 4689.4998
$ws.Range("L132").Value = 151500
$ws.Range("M132").Value = -2159.4998
$ws.Range("N132").Value = -156560

$ws.Range("H135").Value = 1847.0667
$ws.Range("I135").Value = 1943.8462
$ws.Range("J135").Value = 1218
$ws.Range("K135").Value = 17494.6158
$ws.Range("L135").Value = 10962
$ws.Range("M135").Value = -14959.6158
$ws.Range("N135").Value = -16032

$ws.Range("H138").Value = 5513
$ws.Range("I138").Value = 5513
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 16539
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -11399
$ws.Range("N138").ClearContents()

$ws.Range("H141").Value = 3176.1667
$ws.Range("I141").Value = 3188.1738
$ws.Range("J141").Value = 2900
$ws.Range("K141").Value = 9564.5214
$ws.Range("L141").Value = 8700
$ws.Range("M141").Value = -4384.5214
$ws.Range("N141").Value = -19060

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7944289.5
$ws.Range("I61").Value = 12826169
$ws.Range("J61").Value = 11236
$ws.Range("K61").Value = 12826169
$ws.Range("L61").Value = 11236
$ws.Range("M61").Value = -12825957
$ws.Range("N61").Value = -11660

$ws.Range("H74").Value = 6306.241
$ws.Range("I74").Value = 4915.24
$ws.Range("K74").Value = 4915.24
$ws.Range("M74").Value = -4041.24

$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

$ws.Range("H77").Value = 6306.241
$ws.Range("I77").Value = 4915.24
$ws.Range("K77").Value = 24576.2
$ws.Range("M77").Value = -20208.2

$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

$ws.Range("H136").Value = 7944289.5
$ws.Range("I136").Value = 12826169
$ws.Range("J136").Value = 11236
$ws.Range("K136").Value = 38478507
$ws.Range("L136").Value = 33708
$ws.Range("M136").Value = -38475957
$ws.Range("N136").Value = -38808

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3234.8333
$ws.Range("I94").Value = 2720.889
$ws.Range("K94").Value = 2720.889
$ws.Range("M94").Value = -2269.889

$ws.Range("H134").Value = 5328.1665
$ws.Range("I134").Value = 5244
$ws.Range("K134").Value = 15732
$ws.Range("M134").Value = -13197

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2599.2942
$ws.Range("I16").Value = 2094.4546
$ws.Range("J16").Value = 3524.8333
$ws.Range("K16").Value = 2094.4546
$ws.Range("L16").Value = 3524.8333
$ws.Range("M16").Value = -1807.4546
$ws.Range("N16").Value = -4098.8333

$ws.Range("H107").Value = 1370.1428
$ws.Range("I107").Value = 550.6667
$ws.Range("K107").Value = 550.6667
$ws.Range("M107").Value = 1369.3333

$ws.Range("H113").Value = 2599.2942
$ws.Range("I113").Value = 2094.4546
$ws.Range("J113").Value = 3524.8333
$ws.Range("K113").Value = 2094.4546
$ws.Range("L113").Value = 3524.8333
$ws.Range("M113").Value = 75.54539999999997
$ws.Range("N113").Value = -7864.8333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2698.2083
$ws.Range("J68").Value = 3203.5454
$ws.Range("L68").Value = 9610.636200000001
$ws.Range("N68").Value = -11232.6362

$ws.Range("H71").Value = 2698.2083
$ws.Range("J71").Value = 3203.5454
$ws.Range("L71").Value = 28831.9086
$ws.Range("N71").Value = -36943.9086

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 50001
$ws.Range("J74").Value = 50001
$ws.Range("L74").Value = 50001
$ws.Range("N74").Value = -51873

$ws.Range("H77").Value = 50001
$ws.Range("J77").Value = 50001
$ws.Range("L77").Value = 150003
$ws.Range("N77").Value = -159363

$ws.Range("H102").Value = 3235.2
$ws.Range("I102").Value = 2881.8125
$ws.Range("J102").Value = 4648.75
$ws.Range("K102").Value = 2881.8125
$ws.Range("L102").Value = 4648.75
$ws.Range("M102").Value = -1259.8125
$ws.Range("N102").Value = -7892.75

$ws.Range("H107").Value = 442.13333
$ws.Range("I107").Value = 420.625
$ws.Range("J107").Value = 466.7143
$ws.Range("K107").Value = 420.625
$ws.Range("L107").Value = 466.7143
$ws.Range("M107").Value = 1499.375
$ws.Range("N107").Value = -4306.7143

$ws.Range("H113").Value = 667833.3
$ws.Range("I113").Value = 1000750
$ws.Range("K113").Value = 1000750
$ws.Range("M113").Value = -998580

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 6735.9414
$ws.Range("I136").Value = 5403.9165
$ws.Range("J136").Value = 9932.799999999999
$ws.Range("K136").Value = 16211.7495
$ws.Range("L136").Value = 29798.4
$ws.Range("M136").Value = -13661.7495
$ws.Range("N136").Value = -34898.39999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 935.1667
$ws.Range("I100").Value = 935.1667
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1870.3334
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -1329.3334
$ws.Range("N100").ClearContents()

$ws.Range("H132").Value = 6543.34
$ws.Range("I132").Value = 5998.6
$ws.Range("K132").Value = 17995.8
$ws.Range("M132").Value = -15465.8

$ws.Range("H136").Value = 6024.394
$ws.Range("I136").Value = 5510.2
$ws.Range("K136").Value = 16530.6
$ws.Range("M136").Value = -13980.6
